# "added - blog videos"
# Populate the (until now empty) "Videos" worksheet with the IFRAME embed
# table, make it the active/selected sheet, and restore the previously
# selected "Imagenes" sheet's cursor back to A1.

$wb = $excel.ActiveWorkbook
$wsImagenes = $wb.Worksheets.Item(1)
$wsVideos   = $wb.Worksheets.Item(2)

# ---- Videos sheet: column widths ----------------------------------------
$wsVideos.Columns.Item(1).ColumnWidth = 23.07
$wsVideos.Columns.Item(2).ColumnWidth = 94.42

# ---- Header row (A1:B1) ---------------------------------------------------
$wsVideos.Range("A1").Value = "TITULO"
$wsVideos.Range("B1").Value = "IFRAME"

# Re-use the existing dark header formatting (white text / dark fill) that
# is already applied to the "Imagenes" sheet's header row, so the two
# sheets stay visually consistent and no duplicate style is created.
$wsImagenes.Range("A1:B1").Copy()
$wsVideos.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Video rows -------------------------------------------------------
$wsVideos.Range("A2").Value = "Componentes de un exitoso centro comercial | Planigrupo"
$wsVideos.Range("B2").Value = '<iframe width="990" height="743" src="https://www.youtube.com/embed/HIUrz8sKWgI" title="Componentes de un exitoso centro comercial | Planigrupo" frameborder="0" allow="accelerometer; autoplay; clipboard-write; encrypted-media; gyroscope; picture-in-picture; web-share" allowfullscreen></iframe>'

$wsVideos.Range("A3").Value = "La tecnología que transforma los centros comerciales | Planigrupo"
$wsVideos.Range("B3").Value = '<iframe width="990" height="558" src="https://www.youtube.com/embed/vH4mAj65qFg" title="La tecnología que transforma los centros comerciales" frameborder="0" allow="accelerometer; autoplay; clipboard-write; encrypted-media; gyroscope; picture-in-picture; web-share" allowfullscreen></iframe>'

$wsVideos.Range("A4").Value = "Qué considerar al elegir una ubicación para su centro comercial | Planigrupo"
$wsVideos.Range("B4").Value = '<iframe width="990" height="743" src="https://www.youtube.com/embed/Fh_gyPJ-Ib4" title="Qué considerar al elegir una ubicación para su centro comercial | Planigrupo" frameborder="0" allow="accelerometer; autoplay; clipboard-write; encrypted-media; gyroscope; picture-in-picture; web-share" allowfullscreen></iframe>'

# Title cells (column A, rows 2-4) get a larger wrapped font; row 2 is
# vertically centred, rows 3-4 sit on the baseline.
$titleFont = $wsVideos.Range("A2:A4")
$titleFont.Font.Size = 11
$titleFont.WrapText = $true

$wsVideos.Range("A2").VerticalAlignment = -4108
$wsVideos.Range("A3:A4").VerticalAlignment = -4107

# Row heights to fit the wrapped titles.
$wsVideos.Rows.Item(2).RowHeight = 39.55
$wsVideos.Rows.Item(3).RowHeight = 39.55
$wsVideos.Rows.Item(4).RowHeight = 52.2

# ---- Selection / active sheet bookkeeping -------------------------------
$wsImagenes.Range("A1").Select()
$wsVideos.Range("B7").Select()
$wsVideos.Activate()
